$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in column C (Forandrad) from 45184 to 45186 for all data rows (2-172)
$ws.Range("C2:C172").Value = 45186

# Add a friendly display-text second argument to the HYPERLINK formulas in rows 2-9
# Row 2: A 27902-2023
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 27902-2023.xlsx", "A 27902-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 27902-2023.png", "A 27902-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 27902-2023.docx", "A 27902-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 27902-2023.docx", "A 27902-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 27902-2023.docx", "A 27902-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 27902-2023.docx", "A 27902-2023")'

# Row 3: A 8457-2019
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 8457-2019.xlsx", "A 8457-2019")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 8457-2019.png", "A 8457-2019")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 8457-2019.docx", "A 8457-2019")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 8457-2019.docx", "A 8457-2019")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 8457-2019.docx", "A 8457-2019")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 8457-2019.docx", "A 8457-2019")'

# Row 4: A 8446-2019
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 8446-2019.xlsx", "A 8446-2019")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 8446-2019.png", "A 8446-2019")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 8446-2019.docx", "A 8446-2019")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 8446-2019.docx", "A 8446-2019")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 8446-2019.docx", "A 8446-2019")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 8446-2019.docx", "A 8446-2019")'

# Row 5: A 15070-2021
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 15070-2021.xlsx", "A 15070-2021")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 15070-2021.png", "A 15070-2021")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 15070-2021.docx", "A 15070-2021")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 15070-2021.docx", "A 15070-2021")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 15070-2021.docx", "A 15070-2021")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 15070-2021.docx", "A 15070-2021")'

# Row 6: A 44327-2021
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 44327-2021.xlsx", "A 44327-2021")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 44327-2021.png", "A 44327-2021")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 44327-2021.docx", "A 44327-2021")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 44327-2021.docx", "A 44327-2021")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 44327-2021.docx", "A 44327-2021")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 44327-2021.docx", "A 44327-2021")'

# Row 7: A 46919-2021
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 46919-2021.xlsx", "A 46919-2021")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 46919-2021.png", "A 46919-2021")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 46919-2021.docx", "A 46919-2021")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 46919-2021.docx", "A 46919-2021")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 46919-2021.docx", "A 46919-2021")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 46919-2021.docx", "A 46919-2021")'

# Row 8: A 19517-2023
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 19517-2023.xlsx", "A 19517-2023")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 19517-2023.png", "A 19517-2023")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 19517-2023.docx", "A 19517-2023")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 19517-2023.docx", "A 19517-2023")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 19517-2023.docx", "A 19517-2023")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 19517-2023.docx", "A 19517-2023")'

# Row 9: A 24278-2023
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/artfynd/A 24278-2023.xlsx", "A 24278-2023")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/kartor/A 24278-2023.png", "A 24278-2023")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomål/A 24278-2023.docx", "A 24278-2023")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/klagomålsmail/A 24278-2023.docx", "A 24278-2023")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsyn/A 24278-2023.docx", "A 24278-2023")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HOOR/tillsynsmail/A 24278-2023.docx", "A 24278-2023")'

